$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 23:19:11"
$wsZhCn.Range("H2").Value = "2016-03-24 23:19:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 23:19:16"
$wsDeDe.Range("H2").Value = "2016-03-24 23:19:43"
